$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddEmployee")

# Update row 3 - Jane -> Donald, Doe -> Trump
$ws.Range("A3").Value = "Donald"
$ws.Range("C3").Value = "Trump"

# Update row 4 - James -> Katie, Doe -> Ball
$ws.Range("A4").Value = "Katie"
$ws.Range("C4").Value = "Ball"

# Update row 5 - Mark -> Mohammed, Doe -> Salah
$ws.Range("A5").Value = "Mohammed"
$ws.Range("C5").Value = "Salah"

# Add new Employee ID column
$ws.Range("D1").Value = "Employee ID"
$ws.Range("D2").Value = 55555555
$ws.Range("D3").Value = 66666666
$ws.Range("D4").Value = 3333333333
$ws.Range("D5").Value = 4444444

# Column D width auto-fit (bestFit) to match header width
$ws.Columns.Item(4).ColumnWidth = 19.2

# Set the active selection to C11 as in the final file
$ws.Range("C11").Select()
